$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "07_02_2024"
$ws.Range("G2").Value = 7
$ws.Range("G3").Value = 7
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 2.5

[void]$ws.Range("G5").Select()
